# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Sun Jul  7 22:15:45 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.674.01'
$ws.Range("E2").Value = '  -2.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.959.75'
$ws.Range("E3").Value = '  -3.32%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.64'
$ws.Range("E5").Value = '  -5.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.73'
$ws.Range("E6").Value = '  -5.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -5.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.17'
$ws.Range("E9").Value = '  -6.08%  '
$ws.Range("E10").Value = '  -6.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.352'
$ws.Range("E11").Value = '  -4.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.463.45'
$ws.Range("E12").Value = '  -3.44%  '
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("E14").Value = '  -5.19%  '
$ws.Range("E15").Value = '  -8.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.658.03'
$ws.Range("E16").Value = '  -2.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.98'
$ws.Range("E17").Value = '  -3.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.957.14'
$ws.Range("E18").Value = '  -3.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.46'
$ws.Range("E19").Value = '  -5.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.74'
$ws.Range("E20").Value = '  -6.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '316.68'
$ws.Range("E21").Value = '  -7.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("E24").Value = '  -3.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.58'
$ws.Range("E25").Value = '  -3.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0865'
$ws.Range("E28").Value = '  -11.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.47'
$ws.Range("E29").Value = '  -7.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -5.56%  '
$ws.Range("E31").Value = '  -5.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.91'
$ws.Range("E32").Value = '  -5.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.13'
$ws.Range("E33").Value = '  -8.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.98'
$ws.Range("E34").Value = '  -3.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.46'
$ws.Range("E35").Value = '  -7.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.68'
$ws.Range("E36").Value = '  -5.17%  '
$ws.Range("E37").Value = '  -9.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.86'
$ws.Range("E38").Value = '  -9.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0653'
$ws.Range("E39").Value = '  -7.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.989.45'
$ws.Range("E40").Value = '  -3.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.31'
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.640'
$ws.Range("E43").Value = '  -3.93%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.68'
$ws.Range("E44").Value = '  -6.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.155.37'
$ws.Range("E45").Value = '  -8.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.35'
$ws.Range("E46").Value = '  -8.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.84'
$ws.Range("E47").Value = '  -3.47%  '
$ws.Range("E48").Value = '  -11.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0231'
$ws.Range("E49").Value = '  -5.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.96'
$ws.Range("E50").Value = '  -6.43%  '
$ws.Range("E51").Value = '  -13.48%  '
